$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Added CDS All studies testcase": the SamplesTab query (B3) is simplified
# back down to its original column set, dropping the Tumor-status /
# Analyte-Type columns that had been tacked on for a different test case.
$newSamplesQuery = @"
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
   s.phs_accession = 'phs001524' AND gi.reference_genome_assembly = 'GRCh37'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
"@

$ws.Range("B3").Value = $newSamplesQuery

# The saved selection/view moved from B4 up to B3 (row 3 scrolled to top).
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("B3").Select()
